$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.515.50"
$ws.Range("E2").Value = "  +3.83%  "
$ws.Range("D3").Value = "1.738.87"
$ws.Range("E3").Value = "  +4.41%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.44%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4801"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.24%  "
$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2680"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.05%  "
$ws.Range("D9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06242"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").Value = "1.737.88"
$ws.Range("E10").Value = "  +4.22%  "
$ws.Range("D11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07127"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.57%  "
$ws.Range("E12").Value = "  +7.85%  "
$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6210"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.86%  "
$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.541"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.38%  "
$ws.Range("D15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "26.524.76"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006895"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.77%  "
$ws.Range("E20").Value = "  +3.58%  "
$ws.Range("D21").Value = "1.961.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.598"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.07%  "
$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.903"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.35%  "
$ws.Range("D24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.348"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.76%  "
$ws.Range("D27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.814"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.31%  "
$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.425"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.47%  "
$ws.Range("D29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("D30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.995"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.739"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.82%  "
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04587"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.54%  "
$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.616"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.001"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.05%  "
$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6378"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.38%  "
$ws.Range("D37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9285"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "112.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.57%  "
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("E40").Value = "  +8.82%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01515"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.80%  "
$ws.Range("D43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.732"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.52%  "
$ws.Range("D44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3910"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.52%  "
$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.972"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.99%  "
$ws.Range("E46").Value = "  +8.30%  "
$ws.Range("D47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05333"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.52%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.867"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.70%  "
$ws.Range("D50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.253"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.24%  "
$ws.Range("D51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3447"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.40%  "
